# Swap the species-record data held in row 15 and row 16 (columns A, B, D,
# E, F, G, H, Q, R). Everything else on those two rows (C, I, P, S, T, U, V,
# W, Y, Z, AA, AB, AD, AE, AG, AT, AW, AX, AY) is identical between the rows
# and is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R")

foreach ($col in $cols) {
    $cell15 = $ws.Range($col + "15")
    $cell16 = $ws.Range($col + "16")

    $v15 = $cell15.Value2
    $v16 = $cell16.Value2

    $cell15.Value = $v16
    $cell16.Value = $v15
}
